$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (10-14) appended to the "procedure" sheet, covering the
# merged notification / procedure / training-board features.
$rows = @(
    @{ A = 10; B = "Quy trình Offline Lab";      C = "BanNhanSu"; D = "10/09/2022"; E = '<p>Ngày 20/8/2022, Lab thầy Sinh tổ chức offline Lab tại trụ sở Tư Đình, Long Biên, Hà Nội</p>'; I = "Chưa duyệt"; J = "2hfsdjfhjkadfhads" },
    @{ A = 11; B = "Quy trình Offline Lab";      C = "BanDaoTao"; D = "10/09/2022"; E = '<p style="text-align: center;"><span style="color: rgb(255, 0, 0);">Đây là quy trình rất quan trọng</span></p>'; I = "Chưa duyệt"; J = "fjdfhjdfhajf" },
    @{ A = 12; B = "Quy trình Offline PT";       C = "BanDaoTao"; D = "10/09/2022"; E = '<p style="text-align: center;"><span style="background-color: rgb(255, 0, 0);">Đây là quy trình quan trọng</span></p>'; I = "Chưa duyệt"; J = "dfjadfjdsaf" },
    @{ A = 13; B = "Quy trình Offline các PT";   C = "BanDaoTao"; D = "10/09/2022"; E = '<p>đây là quy trình offline PT</p><p><br></p>'; I = "Chưa duyệt"; J = "đâf" },
    @{ A = 14; B = "Nguyễn Văn Bình";            C = "BanDaoTao"; D = "10/09/2022"; E = '<p>Quy trình gặp mặt</p>'; I = "Chưa duyệt"; J = "fjdfhjdfhajf" }
)

foreach ($row in $rows) {
    $r = $row.A
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C

    # Keep column D as a literal text string ("10/09/2022") rather than
    # letting Excel auto-convert it to a date serial value.
    $ws.Cells.Item($r, 4).Value = "'" + $row.D
    $ws.Cells.Item($r, 4).Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $False
    $ws.Cells.Item($r, 7).Value = $False
    $ws.Cells.Item($r, 8).Value = $False
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
}
